$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to round-trip literal text through the clipboard so
# that numeric-looking strings (e.g. "26.09", "1.00") are written back
# verbatim as text instead of being auto-parsed into numbers by Excel's
# normal cell-input heuristics (which a plain `.Value = "26.09"` triggers).
function Set-ExactText($addr, $text) {
    $ws.Range("Z1").Formula = '="' + $text + '"'
    $ws.Range("Z1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-ExactText "D2" '67.759.80'
$ws.Range("E2").Value = "  -1.98%  "

Set-ExactText "D3" '3.756.06'
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  -0.02%  "

Set-ExactText "D5" '597.78'
$ws.Range("E5").Value = "  -2.80%  "

Set-ExactText "D6" '176.56'
$ws.Range("E6").Value = "  -0.17%  "

Set-ExactText "D7" '3.757.85'
$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("E8").Value = "  -0.03%  "

Set-ExactText "D9" '0.528'
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("E10").Value = "  -4.03%  "

$ws.Range("E12").Value = "  -3.75%  "

Set-ExactText "D13" '38.61'
$ws.Range("E13").Value = "  -3.32%  "

Set-ExactText "D14" '0.0000245'
$ws.Range("E14").Value = "  -3.06%  "

Set-ExactText "D15" '4.387.92'
$ws.Range("E15").Value = "  +0.43%  "

Set-ExactText "D16" '3.758.40'
$ws.Range("E16").Value = "  +0.34%  "

Set-ExactText "D17" '67.707.16'
$ws.Range("E17").Value = "  -2.22%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-ExactText "D18" '0.115'
$ws.Range("E18").Value = "  -4.15%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-ExactText "D19" '7.20'
$ws.Range("E19").Value = "  -3.41%  "

Set-ExactText "D20" '16.50'
$ws.Range("E20").Value = "  +1.03%  "

Set-ExactText "D21" '490.51'
$ws.Range("E21").Value = "  -1.56%  "

Set-ExactText "D22" '9.03'
$ws.Range("E22").Value = "  -3.07%  "

Set-ExactText "D23" '0.741'
$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-ExactText "D24" '85.35'
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-ExactText "D25" '0.0000149'
$ws.Range("E25").Value = "  +11.50%  "

Set-ExactText "D26" '2.36'
$ws.Range("E26").Value = "  -6.55%  "

Set-ExactText "D27" '12.28'
$ws.Range("E27").Value = "  -4.08%  "

Set-ExactText "D28" '10.21'
$ws.Range("E28").Value = "  -4.46%  "

Set-ExactText "D29" '1.00'
$ws.Range("E29").Value = "  -0.11%  "

Set-ExactText "D30" '2.95'
$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  -2.82%  "

Set-ExactText "D32" '32.11'
$ws.Range("E32").Value = "  +5.36%  "

Set-ExactText "D33" '7.72'
$ws.Range("E33").Value = "  -3.24%  "

$ws.Range("E34").Value = "  -4.03%  "

Set-ExactText "D35" '1.00'
$ws.Range("E35").Value = "  +0.11%  "

Set-ExactText "D36" '1.00'
$ws.Range("E36").Value = "  -4.01%  "

Set-ExactText "D37" '5.79'
$ws.Range("E37").Value = "  -5.02%  "

$ws.Range("E38").Value = "  -1.70%  "

Set-ExactText "D39" '0.328'
$ws.Range("E39").Value = "  -5.16%  "

Set-ExactText "D40" '446.71'
$ws.Range("E40").Value = "  +0.18%  "

Set-ExactText "D41" '49.13'
$ws.Range("E41").Value = "  -1.08%  "

$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-ExactText "D42" '2.93'
$ws.Range("E42").Value = "  -2.74%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-ExactText "D43" '2.00'
$ws.Range("E43").Value = "  -2.93%  "

Set-ExactText "D44" '8.36'
$ws.Range("E44").Value = "  -2.29%  "

Set-ExactText "D45" '41.20'
$ws.Range("E45").Value = "  -7.84%  "

Set-ExactText "D46" '2.827.01'
$ws.Range("E46").Value = "  -3.91%  "

$ws.Range("E47").Value = "  +0.01%  "

Set-ExactText "D48" '138.78'
$ws.Range("E48").Value = "  +0.19%  "

Set-ExactText "D49" '0.0350'
$ws.Range("E49").Value = "  -2.24%  "

Set-ExactText "D50" '26.09'
$ws.Range("E50").Value = "  -4.28%  "

Set-ExactText "D51" '23.69'
$ws.Range("E51").Value = "  +7.70%  "

$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()